# Natmi following Dr Hou advice
# Updates Ligand/Receptor-expressing cell counts (1 -> 3) and the
# downstream NATMI-derived statistics recomputed for this LR pair sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = [double]"3"
$ws.Range("G2").Value = [double]"20.32821866666667"
$ws.Range("H2").Value = [double]"60.984656"
$ws.Range("I2").Value = [double]"0.004181898474048532"
$ws.Range("J2").Value = [double]"0.004181898474048532"
$ws.Range("K2").Value = [double]"3"
$ws.Range("M2").Value = [double]"6.546185"
$ws.Range("N2").Value = [double]"19.638555"
$ws.Range("O2").Value = [double]"0.06829436374455893"
$ws.Range("P2").Value = [double]"0.06829436374455893"
$ws.Range("Q2").Value = [double]"133.0722801124533"
$ws.Range("R2").Value = [double]"1197.65052101208"
$ws.Range("S2").Value = [double]"0.0002856000955294863"
$ws.Range("T2").Value = [double]"0.0002856000955294863"
$ws.Range("E3").Value = [double]"3"
$ws.Range("G3").Value = [double]"20.32821866666667"
$ws.Range("H3").Value = [double]"60.984656"
$ws.Range("I3").Value = [double]"0.004181898474048532"
$ws.Range("J3").Value = [double]"0.004181898474048532"
$ws.Range("K3").Value = [double]"3"
$ws.Range("M3").Value = [double]"71.03440333333334"
$ws.Range("N3").Value = [double]"213.10321"
$ws.Range("O3").Value = [double]"0.7410803971510699"
$ws.Range("P3").Value = [double]"0.7410803971510698"
$ws.Range("Q3").Value = [double]"1444.002883816196"
$ws.Range("R3").Value = [double]"12996.02595434576"
$ws.Range("S3").Value = [double]"0.003099122981993339"
$ws.Range("T3").Value = [double]"0.003099122981993339"
$ws.Range("E4").Value = [double]"3"
$ws.Range("G4").Value = [double]"20.32821866666667"
$ws.Range("H4").Value = [double]"60.984656"
$ws.Range("I4").Value = [double]"0.004181898474048532"
$ws.Range("J4").Value = [double]"0.004181898474048532"
$ws.Range("K4").Value = [double]"3"
$ws.Range("M4").Value = [double]"0.3068453333333334"
$ws.Range("N4").Value = [double]"0.920536"
$ws.Range("O4").Value = [double]"0.003201224347919758"
$ws.Range("P4").Value = [double]"0.003201224347919758"
$ws.Range("Q4").Value = [double]"6.237619032846223"
$ws.Range("R4").Value = [double]"56.138571295616"
$ws.Range("S4").Value = [double]"1.338719521565264E-05"
$ws.Range("T4").Value = [double]"1.338719521565264E-05"
$ws.Range("E5").Value = [double]"3"
$ws.Range("G5").Value = [double]"20.32821866666667"
$ws.Range("H5").Value = [double]"60.984656"
$ws.Range("I5").Value = [double]"0.004181898474048532"
$ws.Range("J5").Value = [double]"0.004181898474048532"
$ws.Range("K5").Value = [double]"3"
$ws.Range("M5").Value = [double]"17.965059"
$ws.Range("N5").Value = [double]"53.895177"
$ws.Range("O5").Value = [double]"0.1874240147564516"
$ws.Range("P5").Value = [double]"0.1874240147564516"
$ws.Range("Q5").Value = [double]"365.197647711568"
$ws.Range("R5").Value = [double]"3286.778829404112"
$ws.Range("S5").Value = [double]"0.0007837882013100545"
$ws.Range("T5").Value = [double]"0.0007837882013100545"
$ws.Range("E6").Value = [double]"3"
$ws.Range("G6").Value = [double]"4809.896321333334"
$ws.Range("H6").Value = [double]"14429.688964"
$ws.Range("I6").Value = [double]"0.9894865072215304"
$ws.Range("J6").Value = [double]"0.9894865072215304"
$ws.Range("K6").Value = [double]"3"
$ws.Range("M6").Value = [double]"6.546185"
$ws.Range("N6").Value = [double]"19.638555"
$ws.Range("O6").Value = [double]"0.06829436374455893"
$ws.Range("P6").Value = [double]"0.06829436374455893"
$ws.Range("Q6").Value = [double]"31486.47115026745"
$ws.Range("R6").Value = [double]"283378.240352407"
$ws.Range("S6").Value = [double]"0.06757635144452033"
$ws.Range("T6").Value = [double]"0.06757635144452033"
$ws.Range("E7").Value = [double]"3"
$ws.Range("G7").Value = [double]"4809.896321333334"
$ws.Range("H7").Value = [double]"14429.688964"
$ws.Range("I7").Value = [double]"0.9894865072215304"
$ws.Range("J7").Value = [double]"0.9894865072215304"
$ws.Range("K7").Value = [double]"3"
$ws.Range("M7").Value = [double]"71.03440333333334"
$ws.Range("N7").Value = [double]"213.10321"
$ws.Range("O7").Value = [double]"0.7410803971510699"
$ws.Range("P7").Value = [double]"0.7410803971510698"
$ws.Range("Q7").Value = [double]"341668.1152811084"
$ws.Range("R7").Value = [double]"3075013.037529975"
$ws.Range("S7").Value = [double]"0.7332890537473568"
$ws.Range("T7").Value = [double]"0.7332890537473566"
$ws.Range("E8").Value = [double]"3"
$ws.Range("G8").Value = [double]"4809.896321333334"
$ws.Range("H8").Value = [double]"14429.688964"
$ws.Range("I8").Value = [double]"0.9894865072215304"
$ws.Range("J8").Value = [double]"0.9894865072215304"
$ws.Range("K8").Value = [double]"3"
$ws.Range("M8").Value = [double]"0.3068453333333334"
$ws.Range("N8").Value = [double]"0.920536"
$ws.Range("O8").Value = [double]"0.003201224347919758"
$ws.Range("P8").Value = [double]"0.003201224347919758"
$ws.Range("Q8").Value = [double]"1475.894240018301"
$ws.Range("R8").Value = [double]"13283.0481601647"
$ws.Range("S8").Value = [double]"0.003167568298855642"
$ws.Range("T8").Value = [double]"0.003167568298855642"
$ws.Range("E9").Value = [double]"3"
$ws.Range("G9").Value = [double]"4809.896321333334"
$ws.Range("H9").Value = [double]"14429.688964"
$ws.Range("I9").Value = [double]"0.9894865072215304"
$ws.Range("J9").Value = [double]"0.9894865072215304"
$ws.Range("K9").Value = [double]"3"
$ws.Range("M9").Value = [double]"17.965059"
$ws.Range("N9").Value = [double]"53.895177"
$ws.Range("O9").Value = [double]"0.1874240147564516"
$ws.Range("P9").Value = [double]"0.1874240147564516"
$ws.Range("Q9").Value = [double]"86410.0711966363"
$ws.Range("R9").Value = [double]"777690.6407697267"
$ws.Range("S9").Value = [double]"0.1854535337307979"
$ws.Range("T9").Value = [double]"0.1854535337307979"
$ws.Range("E10").Value = [double]"3"
$ws.Range("G10").Value = [double]"2.69506"
$ws.Range("H10").Value = [double]"8.085180000000001"
$ws.Range("I10").Value = [double]"0.000554424737665286"
$ws.Range("J10").Value = [double]"0.000554424737665286"
$ws.Range("K10").Value = [double]"3"
$ws.Range("M10").Value = [double]"6.546185"
$ws.Range("N10").Value = [double]"19.638555"
$ws.Range("O10").Value = [double]"0.06829436374455893"
$ws.Range("P10").Value = [double]"0.06829436374455893"
$ws.Range("Q10").Value = [double]"17.6423613461"
$ws.Range("R10").Value = [double]"158.7812521149"
$ws.Range("S10").Value = [double]"3.786408470309471E-05"
$ws.Range("T10").Value = [double]"3.786408470309471E-05"
$ws.Range("E11").Value = [double]"3"
$ws.Range("G11").Value = [double]"2.69506"
$ws.Range("H11").Value = [double]"8.085180000000001"
$ws.Range("I11").Value = [double]"0.000554424737665286"
$ws.Range("J11").Value = [double]"0.000554424737665286"
$ws.Range("K11").Value = [double]"3"
$ws.Range("M11").Value = [double]"71.03440333333334"
$ws.Range("N11").Value = [double]"213.10321"
$ws.Range("O11").Value = [double]"0.7410803971510699"
$ws.Range("P11").Value = [double]"0.7410803971510698"
$ws.Range("Q11").Value = [double]"191.4419790475334"
$ws.Range("R11").Value = [double]"1722.9778114278"
$ws.Range("S11").Value = [double]"0.0004108733047793679"
$ws.Range("T11").Value = [double]"0.0004108733047793678"
$ws.Range("E12").Value = [double]"3"
$ws.Range("G12").Value = [double]"2.69506"
$ws.Range("H12").Value = [double]"8.085180000000001"
$ws.Range("I12").Value = [double]"0.000554424737665286"
$ws.Range("J12").Value = [double]"0.000554424737665286"
$ws.Range("K12").Value = [double]"3"
$ws.Range("M12").Value = [double]"0.3068453333333334"
$ws.Range("N12").Value = [double]"0.920536"
$ws.Range("O12").Value = [double]"0.003201224347919758"
$ws.Range("P12").Value = [double]"0.003201224347919758"
$ws.Range("Q12").Value = [double]"0.8269665840533335"
$ws.Range("R12").Value = [double]"7.442699256480001"
$ws.Range("S12").Value = [double]"1.774837969303138E-06"
$ws.Range("T12").Value = [double]"1.774837969303138E-06"
$ws.Range("E13").Value = [double]"3"
$ws.Range("G13").Value = [double]"2.69506"
$ws.Range("H13").Value = [double]"8.085180000000001"
$ws.Range("I13").Value = [double]"0.000554424737665286"
$ws.Range("J13").Value = [double]"0.000554424737665286"
$ws.Range("K13").Value = [double]"3"
$ws.Range("M13").Value = [double]"17.965059"
$ws.Range("N13").Value = [double]"53.895177"
$ws.Range("O13").Value = [double]"0.1874240147564516"
$ws.Range("P13").Value = [double]"0.1874240147564516"
$ws.Range("Q13").Value = [double]"48.41691190854"
$ws.Range("R13").Value = [double]"435.7522071768601"
$ws.Range("S13").Value = [double]"0.0001039125102135204"
$ws.Range("T13").Value = [double]"0.0001039125102135204"
$ws.Range("E14").Value = [double]"3"
$ws.Range("G14").Value = [double]"28.08283533333333"
$ws.Range("H14").Value = [double]"84.24850599999999"
$ws.Range("I14").Value = [double]"0.005777169566755752"
$ws.Range("J14").Value = [double]"0.005777169566755752"
$ws.Range("K14").Value = [double]"3"
$ws.Range("M14").Value = [double]"6.546185"
$ws.Range("N14").Value = [double]"19.638555"
$ws.Range("O14").Value = [double]"0.06829436374455893"
$ws.Range("P14").Value = [double]"0.06829436374455893"
$ws.Range("Q14").Value = [double]"183.8354354165367"
$ws.Range("R14").Value = [double]"1654.51891874883"
$ws.Range("S14").Value = [double]"0.0003945481198060132"
$ws.Range("T14").Value = [double]"0.0003945481198060132"
$ws.Range("E15").Value = [double]"3"
$ws.Range("G15").Value = [double]"28.08283533333333"
$ws.Range("H15").Value = [double]"84.24850599999999"
$ws.Range("I15").Value = [double]"0.005777169566755752"
$ws.Range("J15").Value = [double]"0.005777169566755752"
$ws.Range("K15").Value = [double]"3"
$ws.Range("M15").Value = [double]"71.03440333333334"
$ws.Range("N15").Value = [double]"213.10321"
$ws.Range("O15").Value = [double]"0.7410803971510699"
$ws.Range("P15").Value = [double]"0.7410803971510698"
$ws.Range("Q15").Value = [double]"1994.847451811585"
$ws.Range("R15").Value = [double]"17953.62706630426"
$ws.Range("S15").Value = [double]"0.004281347116940427"
$ws.Range("T15").Value = [double]"0.004281347116940427"
$ws.Range("E16").Value = [double]"3"
$ws.Range("G16").Value = [double]"28.08283533333333"
$ws.Range("H16").Value = [double]"84.24850599999999"
$ws.Range("I16").Value = [double]"0.005777169566755752"
$ws.Range("J16").Value = [double]"0.005777169566755752"
$ws.Range("K16").Value = [double]"3"
$ws.Range("M16").Value = [double]"0.3068453333333334"
$ws.Range("N16").Value = [double]"0.920536"
$ws.Range("O16").Value = [double]"0.003201224347919758"
$ws.Range("P16").Value = [double]"0.003201224347919758"
$ws.Range("Q16").Value = [double]"8.617086968801779"
$ws.Range("R16").Value = [double]"77.55378271921599"
$ws.Range("S16").Value = [double]"1.849401587915955E-05"
$ws.Range("T16").Value = [double]"1.849401587915955E-05"
$ws.Range("E17").Value = [double]"3"
$ws.Range("G17").Value = [double]"28.08283533333333"
$ws.Range("H17").Value = [double]"84.24850599999999"
$ws.Range("I17").Value = [double]"0.005777169566755752"
$ws.Range("J17").Value = [double]"0.005777169566755752"
$ws.Range("K17").Value = [double]"3"
$ws.Range("M17").Value = [double]"17.965059"
$ws.Range("N17").Value = [double]"53.895177"
$ws.Range("O17").Value = [double]"0.1874240147564516"
$ws.Range("P17").Value = [double]"0.1874240147564516"
$ws.Range("Q17").Value = [double]"504.509793650618"
$ws.Range("R17").Value = [double]"4540.588142855562"
$ws.Range("S17").Value = [double]"0.001082780314130153"
$ws.Range("T17").Value = [double]"0.001082780314130153"
